# "Add deactivate users task"
#
# The "General Outline" sheet has a "Process" section (starting at A23)
# that lists, as a set of bullets in column B, what the LDAP user import
# process does:
#
#   B24: New Users will be inserted
#   B25: The active status of all users will be updated
#   B26: Site View users have site affiliations updated - only those affiliations in the file will be retained
#   B27: Site View users will have Permissions updated - only those permissions in the file will be retained
#
# This change adds a new leading bullet, "All current users will be
# deactivated", above the existing ones - so the existing four bullets
# shift down one row (24->25, 25->26, 26->27, 27->28) and the new text
# becomes the new B24.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General Outline")

$firstRow = 24
$lastRow = 27

# Shift the existing bullets down one row, working from the bottom up so
# a lower row's old contents aren't clobbered before they've been copied.
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $ws.Range("B" + $r).Copy()
    $ws.Range("B" + ($r + 1)).PasteSpecial(-4104) | Out-Null
}
$excel.CutCopyMode = $false

# Write the new bullet into the now-vacated first row, re-using the
# formatting that was already on that row.
$ws.Range("B" + $firstRow).Value = "All current users will be deactivated"

# Match the refreshed view state: scrolled so row 13 is at the top, with
# B25 (the "New Users will be inserted" bullet, now one row further down)
# selected.
$ws.Select() | Out-Null
$ws.Range("B" + ($firstRow + 1)).Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
